$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '29.415.10'
$c.Style = "Normal"
$ws.Range("E2").Value = '  +0.05%  '
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.850.38'
$c.Style = "Normal"
$ws.Range("E3").Value = '  +0.17%  '
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.000'
$c.Style = "Normal"
$ws.Range("E4").Value = '  +0.09%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '240.61'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.13%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.6287'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.02%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +0.11%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.07650'
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.2908'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -0.75%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '24.81'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +1.34%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '2.175.45'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +17.53%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.07745'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +0.05%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '5.041'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +0.77%  '
$ws.Range("E14").Value = '  +0.37%  '
$ws.Range("E15").Value = '  -1.32%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '83.43'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -0.38%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '6.174'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -0.11%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '29.548.12'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +0.44%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '228.68'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("E20").Value = '  -0.86%  '
$ws.Range("E21").Value = '  +0.14%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '7.457'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -0.15%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '158.04'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +0.51%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '0.1381'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -1.07%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '8.427'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +0.81%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '17.74'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +0.59%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '1.393'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +6.89%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '1.463'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -0.17%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '0.05591'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +0.17%  '
$ws.Range("E32").Value = '  +0.77%  '
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("E34").Value = '  +0.78%  '
$ws.Range("E35").Value = '  -1.95%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '2.592'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +0.25%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.01805'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +0.27%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '1.228.92'
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '2.730'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -1.68%  '
$ws.Range("E40").Value = '  +0.13%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.9082'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +0.10%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +0.15%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '102.38'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +0.59%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '66.03'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -0.06%  '
$ws.Range("E45").Value = '  +0.24%  '
$ws.Range("B46").Value = 'TheSandbox'
$ws.Range("C46").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.4024'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +0.11%  '
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.00000000117'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -4.22%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '9.024'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +0.82%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '1.682'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +0.30%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.1148'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +2.40%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.05702'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -0.09%  '
